# Insert two new weekly price-report rows (Brócoli, Terminal La Palmera de La Serena)
# ahead of the existing row 660, pushing the existing data (rows 660-705) down to
# rows 662-707. The two new rows report the "Primera"/"Segunda" quality prices for
# the week of date-serial 44714.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at 660 and 661; everything currently at 660.. shifts down by 2.
$ws.Range("A660:A661").EntireRow.Insert()

# New row 660 - Calidad "Primera"
$ws.Cells.Item(660, 1).Value = 8
$ws.Cells.Item(660, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(660, 3).Value = "Coquimbo"
$ws.Cells.Item(660, 4).Value = 44714
$ws.Cells.Item(660, 5).Value = 4
$ws.Cells.Item(660, 6).Value = 100112023
$ws.Cells.Item(660, 7).Value = "Brócoli"
$ws.Cells.Item(660, 8).Value = "Sin especificar"
$ws.Cells.Item(660, 9).Value = "Primera"
$ws.Cells.Item(660, 10).Value = 2200
$ws.Cells.Item(660, 11).Value = 800
$ws.Cells.Item(660, 12).Value = 900
$ws.Cells.Item(660, 13).Value = 850
$ws.Cells.Item(660, 14).Value = "`$/unidad"
$ws.Cells.Item(660, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(660, 16).Value = 850
$ws.Cells.Item(660, 17).Value = 1
$ws.Cells.Item(660, 18).Value = "Hortaliza"

# New row 661 - Calidad "Segunda"
$ws.Cells.Item(661, 1).Value = 8
$ws.Cells.Item(661, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(661, 3).Value = "Coquimbo"
$ws.Cells.Item(661, 4).Value = 44714
$ws.Cells.Item(661, 5).Value = 4
$ws.Cells.Item(661, 6).Value = 100112023
$ws.Cells.Item(661, 7).Value = "Brócoli"
$ws.Cells.Item(661, 8).Value = "Sin especificar"
$ws.Cells.Item(661, 9).Value = "Segunda"
$ws.Cells.Item(661, 10).Value = 1320
$ws.Cells.Item(661, 11).Value = 700
$ws.Cells.Item(661, 12).Value = 750
$ws.Cells.Item(661, 13).Value = 725
$ws.Cells.Item(661, 14).Value = "`$/unidad"
$ws.Cells.Item(661, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(661, 16).Value = 725
$ws.Cells.Item(661, 17).Value = 1
$ws.Cells.Item(661, 18).Value = "Hortaliza"
